# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price refresh to the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 116
$ws.Range("H116").Value = 54650.57
$ws.Range("I116").Value = 105308.6
$ws.Range("J116").Value = 8597.817999999999
$ws.Range("K116").Value = 105308.6
$ws.Range("L116").Value = 8597.817999999999
$ws.Range("M116").Value = -101866.6
$ws.Range("N116").Value = -15481.818
# Row 132
$ws.Range("H132").Value = 39588.594
$ws.Range("I132").Value = 50190.81
$ws.Range("J132").Value = 2480.8333
$ws.Range("K132").Value = 150572.43
$ws.Range("L132").Value = 7442.499899999999
$ws.Range("M132").Value = -148042.43
$ws.Range("N132").Value = -12502.4999
# Row 137
$ws.Range("I137").Value = 37501110
$ws.Range("J137").Value = 1799.8462
$ws.Range("K137").Value = 112503330
$ws.Range("L137").Value = 5399.5386
$ws.Range("M137").Value = -112500780
$ws.Range("N137").Value = -10499.5386
# Row 138
$ws.Range("H138").Value = 5929.6
$ws.Range("I138").Value = 6462
$ws.Range("K138").Value = 19386
$ws.Range("M138").Value = -14246

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 4999
$ws.Range("I45").Value = 4999
$ws.Range("K45").Value = 4999
$ws.Range("M45").Value = -4622
# Row 61
$ws.Range("H61").Value = 2952.8667
$ws.Range("I61").Value = 2845.6924
$ws.Range("K61").Value = 2845.6924
$ws.Range("M61").Value = -2633.6924
# Row 74
$ws.Range("H74").Value = 776275.4
$ws.Range("I74").Value = 4760.1035
$ws.Range("K74").Value = 4760.1035
$ws.Range("M74").Value = -3886.1035
# Row 77
$ws.Range("H77").Value = 776275.4
$ws.Range("I77").Value = 4760.1035
$ws.Range("K77").Value = 23800.5175
$ws.Range("M77").Value = -19432.5175
# Row 102
$ws.Range("H102").Value = 1214.4667
$ws.Range("I102").Value = 1365.5
$ws.Range("K102").Value = 1365.5
$ws.Range("M102").Value = 256.5
# Row 122
$ws.Range("H122").Value = 7810.6
$ws.Range("I122").Value = 7632
$ws.Range("J122").Value = 8078.5
$ws.Range("K122").Value = 22896
$ws.Range("L122").Value = 24235.5
$ws.Range("M122").Value = -20446
$ws.Range("N122").Value = -29135.5
# Row 132
$ws.Range("H132").Value = 45461028
$ws.Range("I132").Value = 5600
$ws.Range("J132").Value = 166675500
$ws.Range("K132").Value = 16800
$ws.Range("L132").Value = 500026500
$ws.Range("M132").Value = -14270
$ws.Range("N132").Value = -500031560
# Row 136
$ws.Range("H136").Value = 2952.8667
$ws.Range("I136").Value = 2845.6924
$ws.Range("K136").Value = 8537.0772
$ws.Range("M136").Value = -5987.0772

$ws = $wb.Worksheets.Item("BSM")
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
# Row 134
$ws.Range("H134").Value = 19233584
$ws.Range("I134").Value = 27780418
$ws.Range("J134").Value = 3203.25
$ws.Range("K134").Value = 83341254
$ws.Range("L134").Value = 9609.75
$ws.Range("M134").Value = -83338719
$ws.Range("N134").Value = -14679.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5381.1465
$ws.Range("I31").Value = 3910.6667
$ws.Range("J31").Value = 5794.7188
$ws.Range("K31").Value = 3910.6667
$ws.Range("L31").Value = 5794.7188
$ws.Range("M31").Value = -3615.6667
$ws.Range("N31").Value = -6384.7188
# Row 34
$ws.Range("H34").Value = 5381.1465
$ws.Range("I34").Value = 3910.6667
$ws.Range("J34").Value = 5794.7188
$ws.Range("K34").Value = 3910.6667
$ws.Range("L34").Value = 5794.7188
$ws.Range("M34").Value = -3708.6667
$ws.Range("N34").Value = -6198.7188
# Row 105
$ws.Range("H105").Value = 2803.9443
$ws.Range("I105").Value = 2746.4666
$ws.Range("K105").Value = 2746.4666
$ws.Range("M105").Value = -999.4666000000002
# Row 122
$ws.Range("H122").Value = 13433.523
$ws.Range("I122").Value = 1534.8667
$ws.Range("J122").Value = 43180.168
$ws.Range("K122").Value = 4604.6001
$ws.Range("L122").Value = 129540.504
$ws.Range("M122").Value = -2154.6001
$ws.Range("N122").Value = -134440.504
# Row 132
$ws.Range("H132").Value = 2812.375
$ws.Range("I132").Value = 3015
$ws.Range("K132").Value = 9045
$ws.Range("M132").Value = -6515

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 3248.476
$ws.Range("I68").Value = 1675.6666
$ws.Range("J68").Value = 3877.6
$ws.Range("K68").Value = 5026.9998
$ws.Range("L68").Value = 11632.8
$ws.Range("M68").Value = -4215.9998
$ws.Range("N68").Value = -13254.8
# Row 71
$ws.Range("H71").Value = 3248.476
$ws.Range("I71").Value = 1675.6666
$ws.Range("J71").Value = 3877.6
$ws.Range("K71").Value = 15080.9994
$ws.Range("L71").Value = 34898.4
$ws.Range("M71").Value = -11024.9994
$ws.Range("N71").Value = -43010.4

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 2863.8
$ws.Range("I113").Value = 2499.5
$ws.Range("J113").Value = 4321
$ws.Range("K113").Value = 2499.5
$ws.Range("L113").Value = 4321
$ws.Range("M113").Value = -329.5
$ws.Range("N113").Value = -8661
# Row 122
$ws.Range("H122").Value = 3257.5557
$ws.Range("I122").Value = 4182.091
$ws.Range("K122").Value = 12546.273
$ws.Range("M122").Value = -10096.273
# Row 126
$ws.Range("H126").Value = 2249
$ws.Range("I126").Value = 2165.6667
$ws.Range("J126").Value = 2499
$ws.Range("K126").Value = 6497.000100000001
$ws.Range("L126").Value = 7497
$ws.Range("N126").Value = -12437
$ws.Range("M126").Value = -4027.000100000001
# Row 132
$ws.Range("H132").Value = 9379.375
$ws.Range("I132").Value = 9379.375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 28138.125
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -25608.125
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3305.158
$ws.Range("J40").Value = 4469.25
$ws.Range("L40").Value = 4469.25
$ws.Range("N40").Value = -4741.25
# Row 46
$ws.Range("H46").Value = 2471.4285
$ws.Range("I46").Value = 999
$ws.Range("J46").Value = 2716.8333
$ws.Range("K46").Value = 999
$ws.Range("L46").Value = 2716.8333
$ws.Range("M46").Value = -811
$ws.Range("N46").Value = -3092.8333
# Row 55
$ws.Range("H55").Value = 418.94736
$ws.Range("I55").Value = 516.3333
$ws.Range("J55").Value = 252
$ws.Range("K55").Value = 516.3333
$ws.Range("L55").Value = 252
$ws.Range("M55").Value = -343.3333
$ws.Range("N55").Value = -598
# Row 68
$ws.Range("H68").Value = 5770.7
$ws.Range("I68").Value = 2212.625
$ws.Range("K68").Value = 2212.625
$ws.Range("M68").Value = -1463.625
# Row 71
$ws.Range("H71").Value = 5770.7
$ws.Range("I71").Value = 2212.625
$ws.Range("K71").Value = 11063.125
$ws.Range("M71").Value = -7319.125
# Row 122
$ws.Range("H122").Value = 3454
$ws.Range("I122").Value = 3535.8
$ws.Range("J122").Value = 3249.5
$ws.Range("K122").Value = 10607.4
$ws.Range("L122").Value = 9748.5
$ws.Range("N122").Value = -14648.5
$ws.Range("M122").Value = -8157.400000000001
# Row 132
$ws.Range("H132").Value = 1941.2
$ws.Range("I132").Value = 1951.5
$ws.Range("J132").Value = 1900
$ws.Range("K132").Value = 5854.5
$ws.Range("L132").Value = 5700
$ws.Range("M132").Value = -3324.5
$ws.Range("N132").Value = -10760
# Row 136
$ws.Range("H136").Value = 100003140
$ws.Range("I136").Value = 2899.125
$ws.Range("J136").Value = 500004100
$ws.Range("K136").Value = 8697.375
$ws.Range("L136").Value = 1500012300
$ws.Range("M136").Value = -6147.375
$ws.Range("N136").Value = -1500017400

$ws = $wb.Worksheets.Item("WVR")
# Row 119
$ws.Range("H119").Value = 40000
$ws.Range("J119").Value = 40000
$ws.Range("L119").Value = 40000
$ws.Range("N119").Value = -49676
# Row 122
$ws.Range("H122").Value = 1575.5883
$ws.Range("I122").Value = 1576
$ws.Range("K122").Value = 4728
$ws.Range("M122").Value = -2278
# Row 126
$ws.Range("H126").Value = 5126
$ws.Range("I126").Value = 4647.1665
$ws.Range("J126").Value = 7999
$ws.Range("K126").Value = 13941.4995
$ws.Range("L126").Value = 23997
$ws.Range("M126").Value = -11471.4995
$ws.Range("N126").Value = -28937
# Row 131
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
# Row 132
$ws.Range("H132").Value = 3248.5334
$ws.Range("I132").Value = 3194.5
$ws.Range("J132").Value = 4005
$ws.Range("K132").Value = 9583.5
$ws.Range("L132").Value = 12015
$ws.Range("M132").Value = -7053.5
$ws.Range("N132").Value = -17075
# Row 136
$ws.Range("H136").Value = 1553
$ws.Range("I136").Value = 1596.75
$ws.Range("J136").Value = 1494.6666
$ws.Range("K136").Value = 4790.25
$ws.Range("L136").Value = 4483.9998
$ws.Range("M136").Value = -2240.25
$ws.Range("N136").Value = -9583.9998
